# Apply the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values below are numeric-looking text (e.g. "10.40", "1.00", "0.0000352")
# that must stay TEXT (matching the original inlineStr cells) instead of being
# auto-coerced to numbers (which would drop significant trailing/leading zeros).
# Mark those specific cells as Text format before writing their new value.
$ws.Range("D4,D5,D6,D7,D9,D10,D11,D12,D14,D15,D18,D19,D21,D22,D23,D26,D27,D28,D30,D32,D34,D36,D38,D39,D40,D41,D43,D44,D45,D46,D48,D49,D50,D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.495.31'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '3.944.18'
$ws.Range('E3').Value = '  +4.36%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '471.26'
$ws.Range('E5').Value = '  +7.72%  '
$ws.Range('D6').Value = '147.72'
$ws.Range('E6').Value = '  +3.80%  '
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +10.13%  '
$ws.Range('D11').Value = '0.0000352'
$ws.Range('E11').Value = '  +11.08%  '
$ws.Range('D12').Value = '43.39'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '4.572.14'
$ws.Range('E13').Value = '  +4.46%  '
$ws.Range('D14').Value = '10.40'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '15.12'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').Value = '3.946.99'
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '19.94'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').Value = '1.16'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('D20').Value = '67.655.85'
$ws.Range('E20').Value = '  +1.37%  '
$ws.Range('D21').Value = '435.09'
$ws.Range('E21').Value = '  +4.42%  '
$ws.Range('D22').Value = '3.39'
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('D23').Value = '14.46'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('D26').Value = '38.59'
$ws.Range('E26').Value = '  +3.91%  '
$ws.Range('D27').Value = '9.95'
$ws.Range('E27').Value = '  +4.83%  '
$ws.Range('D28').Value = '10.16'
$ws.Range('E28').Value = '  +4.19%  '
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').Value = '723.25'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').Value = '13.43'
$ws.Range('E32').Value = '  -2.59%  '
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('D34').Value = '42.35'
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0835'
$ws.Range('E35').Value = '  +23.02%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = '57.95'
$ws.Range('E36').Value = '  +2.59%  '
$ws.Range('E37').Value = '  -2.37%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').Value = '5.36'
$ws.Range('E39').Value = '  -4.59%  '
$ws.Range('D40').Value = '0.0476'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').Value = '3.06'
$ws.Range('E41').Value = '  +5.35%  '
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '0.337'
$ws.Range('E43').Value = '  +3.33%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '3.50'
$ws.Range('E45').Value = '  +5.90%  '
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D46').Value = '2.56'
$ws.Range('E46').Value = '  -3.48%  '
$ws.Range('E47').Value = '  +6.55%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '2.82'
$ws.Range('E48').Value = '  +6.85%  '
$ws.Range('D49').Value = '3.28'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').Value = '147.44'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').Value = '2.89'
$ws.Range('E51').Value = '  +2.04%  '
